$d = $word.ActiveDocument

# Locate the (until-now empty) paragraph that immediately follows the
# sentence "...and its relative words are nouns, adjectives, and adverbs."
# -- that is the paragraph being rewritten by this change.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    if ($t -like "*and its relative words are nouns, adjectives, and adverbs.*") {
        $targetPara = $para.Next()
        break
    }
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="78677197" w14:textId="77777777" w:rsidR="00BF565E" w:rsidRPr="003F7833" w:rsidRDefault="00BF565E" w:rsidP="00BF565E">' +
  '<w:pPr>' +
    '<w:pStyle w:val="ListParagraph"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr>' +
    '<w:ind w:firstLineChars="0"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/>' +
      '<w:sz w:val="24"/>' +
      '<w:szCs w:val="24"/>' +
      '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/>' +
      '<w:sz w:val="24"/>' +
      '<w:szCs w:val="24"/>' +
      '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
    '<w:t>The map shows/demonstrates/illutrates</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
      '<w:sz w:val="24"/>' +
      '<w:szCs w:val="24"/>' +
      '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
    '<w:t>&#8230;</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/>' +
      '<w:sz w:val="24"/>' +
      '<w:szCs w:val="24"/>' +
      '<w:lang w:val="en-GB"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve"> </w:t>' +
  '</w:r>' +
  '</w:p>'

$targetPara.Range.InsertXML($xml)
